# =====================================================================
# Adds a new '2022-Q1' sheet (fund-holdings detail) positioned right
# before 'æ»è®¡' (Total), and refreshes 'æ»è®¡' with an extra
# summary row for 2022-Q1.
#
# Implementation notes:
#  - The former 'æ»è®¡' sheet is renamed to '2022-Q1' in place (so it keeps
#    its original sheetId/relationship slot) and then repopulated with the
#    fund-holdings table.
#  - A brand-new 'æ»è®¡' sheet is created by copying the just-renamed
#    sheet (so it inherits the same page/sheet formatting) and is placed
#    immediately after '2022-Q1'; its content is then replaced with the
#    refreshed summary table.
# =====================================================================

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item('总计')
$q4Sheet = $wb.Worksheets.Item('2021-Q4')

# --- Step 1: rename '总计' -> '2022-Q1' (keeps sheetId/rId) ---
$totalSheet.Name = '2022-Q1'
$newQ1 = $wb.Worksheets.Item('2022-Q1')

# --- Step 2: copy it to make the fresh '总计' sheet right after it ---
$newQ1.Copy($null, $newQ1)
$newTotal = $wb.Worksheets.Item('2022-Q1 (2)')
$newTotal.Name = '总计'

# --- Step 3: bring the E:H column header styling + column-A styling from
#             an existing fund-holdings sheet (2021-Q4) onto '2022-Q1' ---
$q4Sheet.Range('B1:H1').Copy()
$newQ1.Range('B1:H1').PasteSpecial(-4122)
$q4Sheet.Range('A2').Copy()
$newQ1.Range('A2:A41').PasteSpecial(-4122)

# --- Step 4: write the '2022-Q1' header row ---
$headers = @('基金代码','基金名称','基金规模','股票总仓位','仓位占比','持有市值(亿元)','仓位排名')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newQ1.Cells.Item(1, $c + 2).Value2 = $headers[$c]
}

# --- Step 5: write the '2022-Q1' fund-holdings rows (A:H, rows 2-41) ---
$data2022Q1 = @(
    ,@('007490', '南方信息创新混合A', '21.90', '91.22', '4.96', '1.0862', 5)
    ,@('160314', '华夏行业混合(LOF)', '21.89', '92.15', '3.64', '0.7968', 8)
    ,@('007355', '汇添富科技创新灵活配置混合A', '16.80', '85.87', '4.42', '0.7426', 6)
    ,@('009683', '汇添富创新增长一年定期开放混合A', '18.11', '71.21', '2.45', '0.4437', 9)
    ,@('012650', '博时半导体主题混合型证券投资基金A', '7.18', '92.40', '5.28', '0.3791', 6)
    ,@('002095', '博时新收益灵活配置混合A', '7.23', '84.41', '4.34', '0.3138', 7)
    ,@('012651', '博时半导体主题混合型证券投资基金C', '3.14', '92.40', '5.28', '0.1658', 6)
    ,@('512330', '南方中证500信息技术指数ETF', '5.09', '99.74', '2.91', '0.1481', 8)
    ,@('007491', '南方信息创新混合C', '2.59', '91.22', '4.96', '0.1285', 5)
    ,@('050018', '博时行业轮动混合', '2.37', '93.93', '5.27', '0.1249', 8)
    ,@('001125', '博时互联网主题灵活配置混合', '5.34', '84.75', '2.29', '0.1223', 10)
    ,@('007356', '汇添富科技创新灵活配置混合C', '2.52', '85.87', '4.42', '0.1114', 6)
    ,@('005062', '博时中证500指数增强A', '6.14', '90.01', '1.70', '0.1044', 3)
    ,@('002096', '博时新收益灵活配置混合C', '2.34', '84.41', '4.34', '0.1016', 7)
    ,@('000522', '华润元大信息传媒科技混合', '1.50', '70.63', '5.92', '0.0888', 6)
    ,@('013345', '富荣信息技术混合A', '1.96', '90.39', '4.36', '0.0855', 6)
    ,@('013339', '创金合信芯片产业股票A', '1.50', '93.43', '4.99', '0.0748', 7)
    ,@('013346', '富荣信息技术混合C', '1.44', '90.39', '4.36', '0.0628', 6)
    ,@('002311', '创金合信中证500指数增强A', '5.72', '92.95', '1.09', '0.0623', 6)
    ,@('012084', '博时睿弘一年定期开放混合型证券投资基金A', '1.78', '91.29', '3.21', '0.0571', 9)
    ,@('009684', '汇添富创新增长一年定期开放混合C', '1.88', '71.21', '2.45', '0.0461', 9)
    ,@('014193', '汇添富中证芯片产业指数增强A', '1.07', '91.67', '3.82', '0.0409', 9)
    ,@('013340', '创金合信芯片产业股票C', '0.61', '93.43', '4.99', '0.0304', 7)
    ,@('002316', '创金合信中证500指数增强C', '2.60', '92.95', '1.09', '0.0283', 6)
    ,@('003241', '创金合信量化发现灵活配置混合A', '2.07', '90.96', '1.07', '0.0221', 7)
    ,@('009608', '广发中证500指数增强A', '1.75', '93.00', '1.13', '0.0198', 7)
    ,@('004930', '华润元大价值优选混合A', '0.32', '65.19', '6.06', '0.0194', 2)
    ,@('005795', '博时中证500指数增强C', '1.14', '90.01', '1.70', '0.0194', 3)
    ,@('007903', '长城量化小盘股票', '1.36', '90.57', '1.17', '0.0159', 3)
    ,@('011377', '创金合信积极成长股票A', '0.29', '94.90', '4.39', '0.0127', 9)
    ,@('014194', '汇添富中证芯片产业指数增强C', '0.29', '91.67', '3.82', '0.0111', 9)
    ,@('004931', '华润元大价值优选混合C', '0.18', '65.19', '6.06', '0.0109', 2)
    ,@('003242', '创金合信量化发现灵活配置混合C', '0.81', '90.96', '1.07', '0.0087', 7)
    ,@('009882', '华润元大核心动力混合A', '0.22', '68.63', '3.96', '0.0087', 8)
    ,@('009609', '广发中证500指数增强C', '0.43', '93.00', '1.13', '0.0049', 7)
    ,@('011378', '创金合信积极成长股票C', '0.11', '94.90', '4.39', '0.0048', 9)
    ,@('009883', '华润元大核心动力混合C', '0.09', '68.63', '3.96', '0.0036', 8)
    ,@('004359', '创金合信量化核心混合A', '0.21', '93.31', '1.61', '0.0034', 10)
    ,@('001607', '英大策略优选混合A', '0.06', '93.31', '1.61', '0.0010', 10)
    ,@('012085', '博时睿弘一年定期开放混合型证券投资基金C', '0.00', '91.29', '3.21', $null, 9)
)

for ($i = 0; $i -lt $data2022Q1.Length; $i++) {
    $r = $i + 2
    $row = $data2022Q1[$i]

    $newQ1.Cells.Item($r, 1).Value2 = $i

    $codeCell = $newQ1.Cells.Item($r, 2)
    $codeCell.Value2 = "'" + $row[0]
    $codeCell.Style = 'Normal'

    $nameCell = $newQ1.Cells.Item($r, 3)
    $nameCell.Value2 = "'" + $row[1]
    $nameCell.Style = 'Normal'

    $sizeCell = $newQ1.Cells.Item($r, 4)
    $sizeCell.Value2 = "'" + $row[2]
    $sizeCell.Style = 'Normal'

    $posCell = $newQ1.Cells.Item($r, 5)
    $posCell.Value2 = "'" + $row[3]
    $posCell.Style = 'Normal'

    $ratioCell = $newQ1.Cells.Item($r, 6)
    $ratioCell.Value2 = "'" + $row[4]
    $ratioCell.Style = 'Normal'

    $mvCell = $newQ1.Cells.Item($r, 7)
    if ($null -eq $row[5]) {
        $mvCell.Value2 = 0
    } else {
        $mvCell.Value2 = "'" + $row[5]
        $mvCell.Style = 'Normal'
    }

    $newQ1.Cells.Item($r, 8).Value2 = $row[6]
}

# --- Step 6: rebuild the '总计' summary sheet (A:D, rows 1-7) ---
$newTotal.Range('A2:A' + 7).Copy()

$dataSummary = @(
    ,@('2022-Q1', 40, 5.51)
    ,@('2021-Q4', 38, 6.16)
    ,@('2021-Q3', 12, 5.37)
    ,@('2021-Q2', 6, 1.51)
    ,@('2021-Q1', 2, 0.2)
    ,@('2020-Q4', 3, 1.08)
)

for ($i = 0; $i -lt $dataSummary.Length; $i++) {
    $r = $i + 2
    $row = $dataSummary[$i]
    $newTotal.Cells.Item($r, 1).Value2 = $i
    $newTotal.Cells.Item($r, 2).Value2 = $row[0]
    $newTotal.Cells.Item($r, 3).Value2 = $row[1]
    $newTotal.Cells.Item($r, 4).Value2 = $row[2]
}

# Make sure column-A styling (s=2) extends to the new 7th row too
$totalSheetForStyle = $wb.Worksheets.Item('2021-Q4')
$totalSheetForStyle.Range('A2').Copy()
$newTotal.Range('A7').PasteSpecial(-4122)

Write-Output 'done'
